$d = $word.ActiveDocument

$replacements = @(
    @("71×24=", "20×18="),
    @("60×63=", "58×46="),
    @("82×54=", "17×13="),
    @("19×22=", "26×17="),
    @("71×62=", "64×46="),
    @("26×65=", "65×86="),
    @("86×80=", "31×90="),
    @("87×88=", "65×84="),
    @("44×13=", "57×60="),
    @("77×60=", "34×66="),
    @("63×15=", "81×48="),
    @("64×57=", "50×90="),
    @("68×11=", "52×11="),
    @("24×91=", "48×45="),
    @("22×81=", "29×89="),
    @("42×27=", "98×83="),
    @("36×86=", "86×81="),
    @("45×26=", "11×83="),
    @("64×33=", "21×73="),
    @("21×54=", "82×35="),
    @("17×42=", "55×30="),
    @("62×68=", "32×60="),
    @("24×45=", "75×33="),
    @("86×45=", "88×55="),
    @("92×13=", "64×87=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
